# Refresh the live cryptos list snapshot (prices + 1h volume deltas).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.968.87'
$ws.Range('E2').Value = '  +1.95%  '
$ws.Range('D3').Value = '2.355.00'
$ws.Range('E3').Value = '  +1.02%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''0.677'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.60%  '
$ws.Range('E6').Value = '  +3.09%  '
$ws.Range('D7').Value = '''72.56'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +11.73%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.537'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +18.93%  '
$ws.Range('D10').Value = '''0.100'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.68%  '
$ws.Range('D11').Value = '''29.39'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +10.71%  '
$ws.Range('E12').Value = '  +2.83%  '
$ws.Range('D13').Value = '2.706.15'
$ws.Range('E13').Value = '  +1.03%  '
$ws.Range('D14').Value = '''16.82'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +10.37%  '
$ws.Range('D15').Value = '''6.67'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.42%  '
$ws.Range('E16').Value = '  +7.96%  '
$ws.Range('D17').Value = '2.367.51'
$ws.Range('E17').Value = '  +1.50%  '
$ws.Range('D18').Value = '43.937.03'
$ws.Range('E18').Value = '  +1.85%  '
$ws.Range('E19').Value = '  +5.17%  '
$ws.Range('D20').Value = '''77.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.13%  '
$ws.Range('E21').Value = '  +4.63%  '
$ws.Range('D22').Value = '''254.85'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.08%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  -3.40%  '
$ws.Range('E25').Value = '  +3.96%  '
$ws.Range('D26').Value = '''10.47'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.73%  '
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('D28').Value = '''22.41'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.47%  '
$ws.Range('D29').Value = '''172.82'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.97%  '
$ws.Range('E30').Value = '  +7.84%  '
$ws.Range('E31').Value = '  +2.49%  '
$ws.Range('E32').Value = '  +5.48%  '
$ws.Range('E33').Value = '  +3.98%  '
$ws.Range('D34').Value = '''0.0723'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.19%  '
$ws.Range('D35').Value = '''5.23'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.66%  '
$ws.Range('D36').Value = '''3.94'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.64%  '
$ws.Range('D37').Value = '''2.44'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.96%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').Value = '''0.0268'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.68%  '
$ws.Range('D40').Value = '''19.46'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.53%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('E42').Value = '  -0.42%  '
$ws.Range('D43').Value = '''1.24'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.90%  '
$ws.Range('D44').Value = '''0.0984'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.74%  '
$ws.Range('E45').Value = '  +1.75%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''98.41'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').Value = '''4.45'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.83%  '
$ws.Range('E48').Value = '  +12.68%  '
$ws.Range('D49').Value = '''2.33'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.82%  '
$ws.Range('D50').Value = '1.438.22'
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('E51').Value = '  +1.30%  '
